$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert two new blank paragraphs at the very top of the document, ahead
#    of the existing "Title Page" paragraph: one will become the new chapter
#    title ("Chapter 6: Creating a Custom Dialog"), the other stays blank
#    (but keeps bold paragraph-mark formatting, matching the target).
# ---------------------------------------------------------------------------
$top = $d.Range(0, 0)
$top.InsertParagraphBefore()
$top.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 2. Fill in the new first paragraph with the chapter title, bold + 16pt
#    (sz/szCs = 32 half-points). A trailing sentinel character is appended
#    temporarily so that the "_GoBack" bookmark can be anchored immediately
#    after the real text but strictly before the paragraph mark (adding a
#    bookmark collapsed at a position flush against a paragraph mark gets
#    mis-anchored, so we park it one character early, then trim the sentinel).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1.Text = "Chapter 6: Creating a Custom DialogX"
$p1.Font.Bold = $true
$p1.Font.Size = 16
$p1.Font.SizeBi = 16

$bmPos = $d.Paragraphs(1).Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the sentinel "X" now that the bookmark is anchored where we want it.
$trimPos = $d.Paragraphs(1).Range.End - 2
$d.Range($trimPos, $trimPos + 1).Delete()

# ---------------------------------------------------------------------------
# 3. The "Title Page" paragraph (now paragraph 3) is untouched.
#    Paragraph 4 ("Hi and wel..come to...") needs its text consolidated into
#    a single run (dropping the stray mid-word bookmark split and the
#    "Bluemix" spell-check proofErr markers). Go through a placeholder value
#    first so the final assignment is recognised as an actual content change.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4).Range
$finalText = "Hi and welcome to the Zero to Cognitive Series. This tutorial is designed to introduce you to building cognitive applications using the IBM Bluemix Cloud. We will use the IBM Cognitive Solution Advisor as our foundation for this tutorial. This is Chapter "

$p4a = $d.Range($p4.Start, $p4.End - 1)
$p4a.Text = "PLACEHOLDER"

$p4b = $d.Range($p4.Start, $d.Paragraphs(4).Range.End - 1)
$p4b.Text = $finalText

Write-Output "Paragraph 1: [$($d.Paragraphs(1).Range.Text)]"
Write-Output "Paragraph 2: [$($d.Paragraphs(2).Range.Text)]"
Write-Output "Paragraph 3: [$($d.Paragraphs(3).Range.Text)]"
Write-Output "Paragraph 4: [$($d.Paragraphs(4).Range.Text)]"
